# Insert two new data rows (124 and 125) into the "Berenjena" price sheet,
# shifting all existing rows 124-204 down to 126-206.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows above the current row 124 (old row 124 becomes row 126, etc.)
$ws.Rows("124:125").Insert()

# --- New row 124 ---
$ws.Range("A124").Value = 9
$ws.Range("B124").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C124").Value = "Metropolitana"
$ws.Range("D124").Value = 44529
$ws.Range("E124").Value = 13
$ws.Range("F124").Value = 100112001
$ws.Range("G124").Value = "Berenjena"
$ws.Range("H124").Value = "Sin especificar"
$ws.Range("I124").Value = "Primera"
$ws.Range("J124").Value = 79
$ws.Range("K124").Value = 8000
$ws.Range("L124").Value = 9000
$ws.Range("M124").Value = 8506
$ws.Range("N124").Value = "$/caja 50 unidades"
$ws.Range("O124").Value = "Región de Arica y Parinacota"
$ws.Range("P124").Value = 170
$ws.Range("Q124").Value = 50
$ws.Range("R124").Value = "Hortaliza"

# --- New row 125 ---
$ws.Range("A125").Value = 9
$ws.Range("B125").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C125").Value = "Metropolitana"
$ws.Range("D125").Value = 44529
$ws.Range("E125").Value = 13
$ws.Range("F125").Value = 100112001
$ws.Range("G125").Value = "Berenjena"
$ws.Range("H125").Value = "Sin especificar"
$ws.Range("I125").Value = "Segunda"
$ws.Range("J125").Value = 43
$ws.Range("K125").Value = 6000
$ws.Range("L125").Value = 6000
$ws.Range("M125").Value = 6000
$ws.Range("N125").Value = "$/caja 100 unidades"
$ws.Range("O125").Value = "Región de Arica y Parinacota"
$ws.Range("P125").Value = 60
$ws.Range("Q125").Value = 100
$ws.Range("R125").Value = "Hortaliza"
